# Commit: "created a model ProdecureItem to load the ul-in UserEngagement-detailsPage"
#
# Semantic changes applied:
#  1. German!O2 ("Ablauf"): periods separating the three sentences become
#     semicolons (in place edit of the existing shared string).
#  2. German!Q2 ("Konkreter_Ablauf"): previously duplicated O2's text verbatim;
#     now gets its own distinct string - same corrected wording, but without
#     the trailing space O2 keeps.
#  3. Row 2 grows taller (409.5 -> 710.2) to fit the longer text.
#  4. The current selection/view on the German sheet moves from A2 to R7.
#  5. The workbook's base font color is pinned to explicit black instead of
#     the theme-1 color.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("German")

# 1) Fix up the "Ablauf" narrative: periods -> semicolons between the three
#    sentences (O2).
$ablaufText = "Con rest voles molor se reptur, erum sum autaquiae prae nonsequat quas ex exero dolupti dolupta tempossimi, volestiure;`n" + `
    "Et fugit od eos eatum expedit, imint as quas comniminus electibustis doloribus. Ent que volore doles eos es dolupta voluptibusam;`n" + `
    "Et fugit od eos eatum expedit, imint as quas comniminus electibustis doloribus. Ent que volore doles eos es dolupta voluptibusam cus sam que nimodios as dolori dolor aut dolorem rem que voluptasimus eum imenihit; "
$ws.Range("O2").Value = $ablaufText

# 2) "Konkreter_Ablauf" (Q2) gets the corrected text too, as its own value
#    (no longer a verbatim duplicate of O2 - trailing space dropped).
$konkreterAblaufText = "Con rest voles molor se reptur, erum sum autaquiae prae nonsequat quas ex exero dolupti dolupta tempossimi, volestiure;`n" + `
    "Et fugit od eos eatum expedit, imint as quas comniminus electibustis doloribus. Ent que volore doles eos es dolupta voluptibusam;`n" + `
    "Et fugit od eos eatum expedit, imint as quas comniminus electibustis doloribus. Ent que volore doles eos es dolupta voluptibusam cus sam que nimodios as dolori dolor aut dolorem rem que voluptasimus eum imenihit;"
$ws.Range("Q2").Value = $konkreterAblaufText

# 3) Row 2 needs to grow to accommodate the (now longer) wrapped text.
$ws.Rows.Item(2).RowHeight = 710.2

# 4) Move the live selection/view to R7.
$ws.Activate()
$ws.Range("R7").Select()

# 5) Pin the base/Normal font color to explicit black (was theme color 1).
$normalStyle = $wb.Styles.Item("Normal")
$normalStyle.Font.Color = 0
